$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 217.05263
$ws.Range("I6").Value = 166.16667
$ws.Range("J6").Value = 304.2857
$ws.Range("K6").Value = 498.50001
$ws.Range("L6").Value = 912.8571000000001
$ws.Range("M6").Value = -386.50001
$ws.Range("N6").Value = -1136.8571
$ws.Range("H8").Value = 650.9
$ws.Range("I8").Value = 1.1666666
$ws.Range("J8").Value = 1625.5
$ws.Range("K8").Value = 3.4999998
$ws.Range("L8").Value = 4876.5
$ws.Range("M8").Value = 135.5000002
$ws.Range("N8").Value = -5154.5
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 5
$ws.Range("K11").Value = 5
$ws.Range("M11").Value = 135
$ws.Range("H41").Value = 1252.875
$ws.Range("I41").Value = 666.7778
$ws.Range("K41").Value = 666.7778
$ws.Range("M41").Value = -226.7778
$ws.Range("H53").Value = 260.125
$ws.Range("I53").Value = 270.8
$ws.Range("J53").Value = 100
$ws.Range("K53").Value = 270.8
$ws.Range("L53").Value = 100
$ws.Range("M53").Value = 366.2
$ws.Range("N53").Value = -1374
$ws.Range("H55").Value = 585
$ws.Range("I55").Value = 188.33333
$ws.Range("K55").Value = 188.33333
$ws.Range("M55").Value = 25.66667000000001
$ws.Range("H80").Value = 1321.7778
$ws.Range("J80").Value = 1319.2
$ws.Range("L80").Value = 3957.6
$ws.Range("N80").Value = -5953.6
$ws.Range("H83").Value = 1321.7778
$ws.Range("J83").Value = 1319.2
$ws.Range("L83").Value = 11872.8
$ws.Range("N83").Value = -21856.8
$ws.Range("H87").Value = 95021.664
$ws.Range("J87").Value = 95021.664
$ws.Range("L87").Value = 95021.664
$ws.Range("N87").Value = -97517.664
$ws.Range("H90").Value = 95021.664
$ws.Range("J90").Value = 95021.664
$ws.Range("L90").Value = 285064.992
$ws.Range("N90").Value = -297544.992
$ws.Range("N100").ClearContents()
$ws.Range("H100").Value = 1143
$ws.Range("I100").Value = 1143
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1143
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -602
$ws.Range("H133").Value = 112500
$ws.Range("J133").Value = 112500
$ws.Range("L133").Value = 112500
$ws.Range("N133").Value = -122620
$ws.Range("H138").Value = 3137.8
$ws.Range("J138").Value = 3089.923
$ws.Range("L138").Value = 9269.769
$ws.Range("N138").Value = -19549.769
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4411.222
$ws.Range("I74").Value = 3281.1667
$ws.Range("K74").Value = 3281.1667
$ws.Range("M74").Value = -2407.1667
$ws.Range("H77").Value = 4411.222
$ws.Range("I77").Value = 3281.1667
$ws.Range("K77").Value = 16405.8335
$ws.Range("M77").Value = -12037.8335
$ws.Range("H110").Value = 847.1539
$ws.Range("I110").Value = 865.7
$ws.Range("J110").Value = 785.3333
$ws.Range("K110").Value = 865.7
$ws.Range("L110").Value = 785.3333
$ws.Range("M110").Value = 1179.3
$ws.Range("N110").Value = -4875.3333

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2945.6
$ws.Range("I105").Value = 2819.9375
$ws.Range("K105").Value = 2819.9375
$ws.Range("M105").Value = -1072.9375

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5431.365
$ws.Range("I31").Value = 2322.8
$ws.Range("J31").Value = 5762.064
$ws.Range("K31").Value = 2322.8
$ws.Range("L31").Value = 5762.064
$ws.Range("M31").Value = -2027.8
$ws.Range("N31").Value = -6352.064
$ws.Range("H34").Value = 5431.365
$ws.Range("I34").Value = 2322.8
$ws.Range("J34").Value = 5762.064
$ws.Range("K34").Value = 2322.8
$ws.Range("L34").Value = 5762.064
$ws.Range("M34").Value = -2120.8
$ws.Range("N34").Value = -6166.064
$ws.Range("N48").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("H122").Value = 1626.9231
$ws.Range("I122").Value = 1365
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4095
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1645
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 4656.643
$ws.Range("I132").Value = 4465.8887
$ws.Range("K132").Value = 13397.6661
$ws.Range("M132").Value = -10867.6661

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 174159.17
$ws.Range("M11").ClearContents()
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("H18").Value = 3023.8125
$ws.Range("I18").Value = 1135.4286
$ws.Range("J18").Value = 4492.5557
$ws.Range("K18").Value = 3406.2858
$ws.Range("L18").Value = 13477.6671
$ws.Range("M18").Value = -3237.2858
$ws.Range("N18").Value = -13815.6671
$ws.Range("H68").Value = 2749.3333
$ws.Range("J68").Value = 2749.75
$ws.Range("L68").Value = 8249.25
$ws.Range("N68").Value = -9871.25
$ws.Range("H71").Value = 2749.3333
$ws.Range("J71").Value = 2749.75
$ws.Range("L71").Value = 24747.75
$ws.Range("N71").Value = -32859.75
$ws.Range("H121").Value = 228.66667
$ws.Range("J121").Value = 300
$ws.Range("L121").Value = 900
$ws.Range("N121").Value = -3520
$ws.Range("H138").Value = 7660
$ws.Range("J138").Value = 8825
$ws.Range("L138").Value = 26475
$ws.Range("N138").Value = -36755
$ws.Range("H139").Value = 2970.6
$ws.Range("I139").Value = 2121.5833
$ws.Range("K139").Value = 6364.749899999999
$ws.Range("M139").Value = -1224.749899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N97").ClearContents()
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1000
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -504
$ws.Range("H99").Value = 3300
$ws.Range("I99").Value = 3300
$ws.Range("K99").Value = 3300
$ws.Range("M99").Value = -1054
$ws.Range("H113").Value = 3982.5
$ws.Range("I113").Value = 2248.75
$ws.Range("K113").Value = 2248.75
$ws.Range("M113").Value = -78.75
$ws.Range("H122").Value = 1974.1538
$ws.Range("I122").Value = 1872.8572
$ws.Range("J122").Value = 2092.3333
$ws.Range("K122").Value = 5618.571599999999
$ws.Range("L122").Value = 6276.999899999999
$ws.Range("M122").Value = -3168.571599999999
$ws.Range("N122").Value = -11176.9999
$ws.Range("H126").Value = 3834
$ws.Range("I126").Value = 3251
$ws.Range("K126").Value = 9753
$ws.Range("M126").Value = -7283

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 819.3333
$ws.Range("I30").Value = 583.2
$ws.Range("K30").Value = 583.2
$ws.Range("M30").Value = -475.2
$ws.Range("H46").Value = 5217.5713
$ws.Range("I46").Value = 3841.8333
$ws.Range("K46").Value = 3841.8333
$ws.Range("M46").Value = -3653.8333
$ws.Range("H55").Value = 1742.75
$ws.Range("I55").Value = 1742.75
$ws.Range("K55").Value = 1742.75
$ws.Range("M55").Value = -1569.75
$ws.Range("H68").Value = 6599.5
$ws.Range("I68").Value = 2998.3333
$ws.Range("K68").Value = 2998.3333
$ws.Range("M68").Value = -2249.3333
$ws.Range("H71").Value = 6599.5
$ws.Range("I71").Value = 2998.3333
$ws.Range("K71").Value = 14991.6665
$ws.Range("M71").Value = -11247.6665
$ws.Range("H130").Value = 49995
$ws.Range("J130").Value = 49995
$ws.Range("L130").Value = 49995
$ws.Range("N130").Value = -60035

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M29").ClearContents()
$ws.Range("H29").Value = 999
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("H86").Value = 70022.664
$ws.Range("J86").Value = 70022.664
$ws.Range("L86").Value = 70022.664
$ws.Range("N86").Value = -72268.664
$ws.Range("H89").Value = 70022.664
$ws.Range("J89").Value = 70022.664
$ws.Range("L89").Value = 350113.32
$ws.Range("N89").Value = -361345.32
